$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1001, 18, 30, 75, 60, 72)
    3  = @(1202, 2, 10, 10, 10, 10)
    4  = @(901, 16, 15, 45, 60, 60)
    5  = @(601, 9, 60, 67, 60, 42)
    6  = @(1201, 2, 10, 10, 10, 10)
    7  = @(101, 9, 30, 15, 60, 15)
    8  = @(902, 1, 0, 0, 0, 0)
    11 = @(1203, 3, 15, 15, 15, 15)
    12 = @(301, 6, 45, 30, 60, 45)
    13 = @(701, 3, 90, 45, 97, 15)
    14 = @(201, 9, 30, 15, 45, 30)
    15 = @(801, 3, 67, 65, 52, 45)
    16 = @(3, 0, 3, 3, 3, 3)
    17 = @(502, 0, 4, 0, 0, 0)
    18 = @(1101, 0, 15, 30, 30, 0)
    19 = @(2, 0, 2, 2, 2, 2)
    20 = @(802, 0, 4, 5, 4, 0)
    21 = @(1, 0, 2, 2, 2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
